$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - set values then copy the header style from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for I and J columns (rows 2-22)
$values = @(
    @(1, 5),
    @(2, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(3, 7),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 8),
    @(1, 5),
    @(1, 6),
    @(3, 7),
    @(1, 3),
    @(1, 5),
    @(1, 8),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 2),
    @(1, 2)
)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
